$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume(1h) figures.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the original inline-string cells) instead of
# auto-converting numeric-looking text to a number; the style is then
# reset to "Normal" so no stray number-format/quote-prefix style sticks
# to the cell.
$c = $ws.Range("D2")
$c.Value = "'315.61"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'39.55"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'2.53%"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.141"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'0.70%"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'0.08176"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'0.76%"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'1.964"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'-0.38%"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'8.227"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'3.56%"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.9266"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'-0.27%"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.1412"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'-1.30%"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.1985"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'1.37%"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.09057"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'-0.19%"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.03507"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'0.11%"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'0.05%"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'0.001392"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'-0.92%"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'0.005943"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'-3.08%"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'3.654"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'-2.11%"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'4.237"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'1.28%"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'-8.96%"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'0.02%"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'0.1303"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'0.69%"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'4.761"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'-0.65%"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'-1.08%"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'0.04379"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'0.58%"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'-0.04%"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'0.004788"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'-0.97%"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'-0.05%"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'0.0003996"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'-10.15%"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.02206"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'5.66%"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.05184"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'1.24%"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.007528"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'0.72%"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.009752"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'-3.75%"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'0.1373"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'1.09%"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.002130"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'-0.05%"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.009126"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'-1.56%"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.00006405"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'2.35%"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'0.00000000749"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'-0.24%"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'-8.70%"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.001199"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'-25.13%"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.00002098"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'-0.24%"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'0.0001998"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'-0.24%"
$c.Style = "Normal"

Write-Output "Updated 70 cells"
